$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "91 x 62" + [char]11 + "  6    2" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "1|    |"
$t.Cell(1,2).Range.Text = "98 x 47" + [char]11 + "  4    7" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "8|    |"
$t.Cell(1,3).Range.Text = "16 x 71" + [char]11 + "  7    1" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "6|    |"

$t.Cell(2,1).Range.Text = "44 x 44" + [char]11 + "  4    4" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "4|    |"
$t.Cell(2,2).Range.Text = "94 x 34" + [char]11 + "  3    4" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "4|    |"
$t.Cell(2,3).Range.Text = "22 x 43" + [char]11 + "  4    3" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "2|    |"

$t.Cell(3,1).Range.Text = "97 x 22" + [char]11 + "  2    2" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "7|    |"
$t.Cell(3,2).Range.Text = "64 x 25" + [char]11 + "  2    5" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "4|    |"
$t.Cell(3,3).Range.Text = "11 x 10" + [char]11 + "  1    0" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "1|    |"

$t.Cell(4,1).Range.Text = "37 x 56" + [char]11 + "  5    6" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "7|    |"
$t.Cell(4,2).Range.Text = "98 x 52" + [char]11 + "  5    2" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "8|    |"
$t.Cell(4,3).Range.Text = "89 x 70" + [char]11 + "  7    0" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "9|    |"

$t.Cell(5,1).Range.Text = "74 x 39" + [char]11 + "  3    9" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "4|    |"
$t.Cell(5,2).Range.Text = "34 x 68" + [char]11 + "  6    8" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "4|    |"
$t.Cell(5,3).Range.Text = "37 x 48" + [char]11 + "  4    8" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "7|    |"
